$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: CheckIn
$ws.Cells.Item(2,1).Value = "CheckIn"
$ws.Cells.Item(2,2).Value = "[0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0]"
$ws.Cells.Item(2,3).Value = 6
$ws.Cells.Item(2,4).Value = '["Business", "Speakers", "Logistics"]'

# Row 3: Auditorio
$ws.Cells.Item(3,1).Value = "Auditorio"
$ws.Cells.Item(3,2).Value = "[0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0]"
$ws.Cells.Item(3,3).Value = 2
$ws.Cells.Item(3,4).Value = '["Logistics"]'

# Row 4: Refeicoes (was Almocos)
$ws.Cells.Item(4,1).Value = "Refeicoes"
$ws.Cells.Item(4,2).Value = "[0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,0,0,0,0,0,1,1,0,0,0,0,0,1,1,0,0,0,0,0,1,1,0,0,0,0,0,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0]"
$ws.Cells.Item(4,3).Value = 6
$ws.Cells.Item(4,4).Value = '["Business", "Logistics"]'

# Row 5: Divulgacao
$ws.Cells.Item(5,1).Value = "Divulgacao"
$ws.Cells.Item(5,2).Value = "[0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0]"
$ws.Cells.Item(5,3).Value = 5
$ws.Cells.Item(5,4).Value = '["Marketing", "Volunteer"]'

# Row 6: Workshops
$ws.Cells.Item(6,1).Value = "Workshops"
$ws.Cells.Item(6,2).Value = "[0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0]"
$ws.Cells.Item(6,3).Value = 4
$ws.Cells.Item(6,4).Value = "[]"

# Row 7: MontagemDesmontagem (new)
$ws.Cells.Item(7,1).Value = "MontagemDesmontagem"
$ws.Cells.Item(7,2).Value = "[1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1]"
$ws.Cells.Item(7,3).Value = 6
$ws.Cells.Item(7,4).Value = "[]"

# Row 8: CoffeeBreak (new)
$ws.Cells.Item(8,1).Value = "CoffeeBreak"
$ws.Cells.Item(8,2).Value = "[0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0]"
$ws.Cells.Item(8,3).Value = 3
$ws.Cells.Item(8,4).Value = "[]"

# Column A width widened to fit the new, longer shift-name strings
$ws.Columns.Item(1).ColumnWidth = 21.166666666666668

# Selection moves to H6
$ws.Range("H6").Select() | Out-Null
